# Add a "Salary" column (K) to the student/employee roster, with per-row
# salary figures formatted as whole-number currency-style numbers (#,##0).
# A handful of rows (no hire-date/degree info on the source sheet) are left
# without a salary figure, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K1").Value = "Salary"

# row number -> salary amount
$salary = @{
    2  = 45000
    3  = 42000
    4  = 38000
    5  = 55000
    6  = 50000
    7  = 55000
    11 = 44000
    12 = 40000
    13 = 55000
    14 = 65000
    16 = 48000
    17 = 44000
    18 = 48000
    19 = 50000
    20 = 35000
    21 = 40000
    22 = 45000
    23 = 38000
    24 = 35000
    25 = 55000
    26 = 55000
    27 = 60000
    29 = 55000
    30 = 50000
    31 = 43000
    33 = 52000
    34 = 65000
    35 = 60000
    36 = 60000
    38 = 50000
    39 = 70000
    40 = 55000
    41 = 48000
    42 = 60000
    43 = 55000
    44 = 50000
    45 = 75000
    47 = 50000
    48 = 48000
    49 = 44000
    50 = 52500
}

foreach ($row in $salary.Keys) {
    $cell = $ws.Range("K" + $row)
    $cell.Value = $salary[$row]
    $cell.NumberFormat = "#,##0"
}

# Mirror the workbook's saved UI state: K1 is the active/selected cell, with
# the grid scrolled so column C leads the visible area.
$ws.Range("K1").Select()
$excel.ActiveWindow.ScrollColumn = 3
